$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were reported between the existing
# 2022-08-?? ("Región de Ñuble") entries and the rest of the series.
# Insert two blank rows right before the current row 223 so every
# subsequent record (old rows 223-261) shifts down to 225-263, then
# populate the two freshly inserted rows with the new data.
$ws.Rows("223:224").Insert()

$ws.Range("A223").Value = 11
$ws.Range("B223").Value = "Vega Monumental Concepción"
$ws.Range("C223").Value = "Bíobío"
$ws.Range("D223").Value = 44798
$ws.Range("E223").Value = 8
$ws.Range("F223").Value = 100114013
$ws.Range("G223").Value = "Zanahoria"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 1000
$ws.Range("K223").Value = 8000
$ws.Range("L223").Value = 8500
$ws.Range("M223").Value = 8250
$ws.Range("N223").Value = "$/saco 20 kilos"
$ws.Range("O223").Value = "Región de La Araucanía"
$ws.Range("P223").Value = 412
$ws.Range("Q223").Value = 20
$ws.Range("R223").Value = "Hortaliza"

$ws.Range("A224").Value = 11
$ws.Range("B224").Value = "Vega Monumental Concepción"
$ws.Range("C224").Value = "Bíobío"
$ws.Range("D224").Value = 44798
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = 100114013
$ws.Range("G224").Value = "Zanahoria"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Segunda"
$ws.Range("J224").Value = 500
$ws.Range("K224").Value = 7000
$ws.Range("L224").Value = 7000
$ws.Range("M224").Value = 7000
$ws.Range("N224").Value = "$/saco 20 kilos"
$ws.Range("O224").Value = "Región de La Araucanía"
$ws.Range("P224").Value = 350
$ws.Range("Q224").Value = 20
$ws.Range("R224").Value = "Hortaliza"
